$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: RandomForestRegressor - values change
$ws.Range("B3").Value = 7245554834597.273
$ws.Range("C3").Value = 6391439790699.064
$ws.Range("D3").Value = 8386556660167.938

# Row 4: GradientBoostingRegressor -> DecisionTreeRegressor
$ws.Range("A4").Value = "DecisionTreeRegressor"
$ws.Range("B4").Value = 1801126611537.544
$ws.Range("C4").Value = 1722816758862.001
$ws.Range("D4").Value = 1801126611537.545

# Row 5: AdaBoostRegressor -> MLPRegressor
$ws.Range("A5").Value = "MLPRegressor"
$ws.Range("B5").Value = 168184780109034.4
$ws.Range("C5").Value = 137570117159743.7
$ws.Range("D5").Value = 212496469990763.2
